$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from an existing header cell (C1) to the new headers (F1:H1)
$ws.Range("C1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# New header cells F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Fill F2:H18 with boolean FALSE
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
